$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to text so values like '1.000' or '316.74'
# are not auto-converted to numbers by Excel's input parser.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Update Price (D) column
$ws.Range("D2").Value = "24.600.08"
$ws.Range("D3").Value = "1.696.63"
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").Value = "316.74"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D7").Value = "0.3946"
$ws.Range("D8").Value = "0.4022"
$ws.Range("D9").Value = "1.532"
$ws.Range("D10").Value = "1.000"
$ws.Range("D11").Value = "52.75"
$ws.Range("D12").Value = "0.08769"
$ws.Range("D13").Value = "7.227"
$ws.Range("D14").Value = "23.31"
$ws.Range("D15").Value = "8.054"
$ws.Range("D16").Value = "0.00001318"
$ws.Range("D17").Value = "1.693.33"
$ws.Range("D18").Value = "99.96"
$ws.Range("D19").Value = "0.07076"
$ws.Range("D20").Value = "19.72"
$ws.Range("D21").Value = "6.936"
$ws.Range("D22").Value = "0.9993"
$ws.Range("D23").Value = "14.21"
$ws.Range("D24").Value = "24.599.06"
$ws.Range("D25").Value = "3.144"
$ws.Range("D27").Value = "23.05"
$ws.Range("D28").Value = "162.62"
$ws.Range("D29").Value = "137.17"
$ws.Range("D30").Value = "5.191"
$ws.Range("D31").Value = "7.519"
$ws.Range("D32").Value = "1.880.02"
$ws.Range("D33").Value = "1.093"
$ws.Range("D34").Value = "0.08589"
$ws.Range("D35").Value = "7.227"
$ws.Range("D36").Value = "11.36"
$ws.Range("D37").Value = "0.2740"
$ws.Range("D38").Value = "1.929"
$ws.Range("D40").Value = "0.09125"
$ws.Range("D42").Value = "1.479"
$ws.Range("D43").Value = "0.7688"
$ws.Range("D44").Value = "0.7197"
$ws.Range("D45").Value = "15.77"
$ws.Range("D46").Value = "2.554"
$ws.Range("D47").Value = "4.220"
$ws.Range("D48").Value = "0.9992"
$ws.Range("D49").Value = "1.331"
$ws.Range("D50").Value = "141.02"
$ws.Range("D51").Value = "0.07996"

# Restore original (default) formatting now that the text values are committed
$priceRange.ClearFormats()

# Update Volume(1h) (E) column
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  +4.75%  "
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  +7.55%  "
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("E15").Value = "  +10.02%  "
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("E20").Value = "  +3.38%  "
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("E25").Value = "  +10.55%  "
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +5.97%  "
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +10.52%  "
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("E35").Value = "  +9.74%  "
$ws.Range("E36").Value = "  +8.43%  "
$ws.Range("E37").Value = "  +3.83%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("E40").Value = "  +3.78%  "
$ws.Range("E41").Value = "  +9.42%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("E49").Value = "  +8.65%  "
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").Value = "  +2.61%  "
